# Generate Report for Handback
# Refresh the handoff/handback datetimes for the 9675f9d4... row (row 2)
# on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-25 07:54:50"
$wsZhCn.Range("G2").Value = "2016-01-25 07:55:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-25 07:55:06"
$wsDeDe.Range("G2").Value = "2016-01-25 07:55:57"
